$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45553, 5270),
    @(45553.01041666666, 5240),
    @(45553.02083333334, 5210),
    @(45553.03125, 5170),
    @(45553.04166666666, 5130),
    @(45553.05208333334, 5100),
    @(45553.0625, 5080),
    @(45553.07291666666, 5070),
    @(45553.08333333334, 5060),
    @(45553.09375, 5050),
    @(45553.10416666666, 5050),
    @(45553.11458333334, 5070),
    @(45553.125, 5100),
    @(45553.13541666666, 5110),
    @(45553.14583333334, 5120),
    @(45553.15625, 5150),
    @(45553.16666666666, 5190),
    @(45553.17708333334, 5250),
    @(45553.1875, 5330),
    @(45553.19791666666, 5430),
    @(45553.20833333334, 5550),
    @(45553.21875, 5680),
    @(45553.22916666666, 5820),
    @(45553.23958333334, 5960),
    @(45553.25, 6090),
    @(45553.26041666666, 6220),
    @(45553.27083333334, 6360),
    @(45553.28125, 6440),
    @(45553.29166666666, 6520),
    @(45553.30208333334, 6550),
    @(45553.3125, 6550),
    @(45553.32291666666, 6550),
    @(45553.33333333334, 6550),
    @(45553.34375, 6540),
    @(45553.35416666666, 6520),
    @(45553.36458333334, 6480),
    @(45553.375, 6430),
    @(45553.38541666666, 6370),
    @(45553.39583333334, 6310),
    @(45553.40625, 6260),
    @(45553.41666666666, 6210),
    @(45553.42708333334, 6180),
    @(45553.4375, 6150),
    @(45553.44791666666, 6140),
    @(45553.45833333334, 6140),
    @(45553.46875, 6150),
    @(45553.47916666666, 6160),
    @(45553.48958333334, 6180),
    @(45553.5, 6200),
    @(45553.51041666666, 6220),
    @(45553.52083333334, 6230),
    @(45553.53125, 6250),
    @(45553.54166666666, 6260),
    @(45553.55208333334, 6270),
    @(45553.5625, 6290),
    @(45553.57291666666, 6310),
    @(45553.58333333334, 6330),
    @(45553.59375, 6360),
    @(45553.60416666666, 6390),
    @(45553.61458333334, 6430),
    @(45553.625, 6460),
    @(45553.63541666666, 6500),
    @(45553.64583333334, 6530),
    @(45553.65625, 6560),
    @(45553.66666666666, 6590),
    @(45553.67708333334, 6610),
    @(45553.6875, 6630),
    @(45553.69791666666, 6650),
    @(45553.70833333334, 6680),
    @(45553.71875, 6710),
    @(45553.72916666666, 6750),
    @(45553.73958333334, 6810),
    @(45553.75, 6880),
    @(45553.76041666666, 6950),
    @(45553.77083333334, 7020),
    @(45553.78125, 7100),
    @(45553.79166666666, 7170),
    @(45553.80208333334, 7200),
    @(45553.8125, 7200),
    @(45553.82291666666, 7180),
    @(45553.83333333334, 7070),
    @(45553.84375, 6940),
    @(45553.85416666666, 6800),
    @(45553.86458333334, 6660),
    @(45553.875, 6500),
    @(45553.88541666666, 6330),
    @(45553.89583333334, 6160),
    @(45553.90625, 6040),
    @(45553.91666666666, 5900),
    @(45553.92708333334, 5790),
    @(45553.9375, 5690),
    @(45553.94791666666, 5600),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $data[$i][1]
}
